$wb = $excel.ActiveWorkbook

# "Presupuesto" sheet: clear the sample/demo values from row 2 (A2:C2),
# leaving the cell styles/formatting intact (D2 was already blank).
$wsPresupuesto = $wb.Worksheets.Item("Presupuesto")
$wsPresupuesto.Range("A2:C2").ClearContents()

# "Priorización" sheet: remove the sample/demo data rows (3 and 5-13),
# leaving only the header rows (1-2) so the template ships empty.
$wsPriorizacion = $wb.Worksheets.Item("Priorización")
$wsPriorizacion.Range("A3:Q13").Clear()

Write-Output "done"
